$d = $word.ActiveDocument

# 1. Update the heading text
$d.Content.Find.Execute(
    "Implante de Marcapasso Convencional com Monitoramento Remoto (Azure™)",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Marcapasso Dupla Câmara (Azure – monitoramento remoto)", 2)

# 2. Remove the now-redundant subtitle paragraph entirely
$removed = $false
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs.Item($i)
    if ($p.Range.Text.TrimEnd([char]13, [char]7) -eq "Marcapasso de câmara dupla com monitoramento remoto.") {
        $p.Range.Delete()
        $removed = $true
        break
    }
}

# 3. Update the materials bullet lines
$d.Content.Find.Execute(
    "Gerador – Azure™",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "• Gerador Azure – marcapasso dupla câmara com monitoramento remoto", 2)

$d.Content.Find.Execute(
    "Eletrodo Ventricular – 5076-52",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "• Eletrodo 5076-52", 2)

$d.Content.Find.Execute(
    "Eletrodo Atrial – 5076-58",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "• Eletrodo 5076-58", 2)

$d.Content.Find.Execute(
    "Introdutor – 2",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "• Introdutor – 2", 2)
